# Celerio role.xlsx regeneration added a new "search_full_text" search
# field to the "Search" sheet, just above the existing roleName search
# row. This inserts a new row 4 (pushing the former row 4 down to row 5)
# and fills it with the new label/value pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# Insert a new blank row above the current row 4 (role_roleName / roleName),
# shifting it down to row 5.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new search_full_text label/value pair.
$ws.Range("A4").Value = "`${msg.getProperty('search_full_text')}"
$ws.Range("B4").Value = "`${search_full_text}"
